$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 13 (shifts existing rows 13-31 down to 15-33)
$ws.Rows("13:14").Insert()

# Keep the table ("Tabla13") in sync with the new used range
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:H33"))

# New row 13: "2.1." Registro asistencia Capacitaciones y Seguimientos (Carpeta)
$ws.Range("A13").Value = "2.1."
$ws.Range("B13").Value = "2."
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = "Registro asistencia Capacitaciones y Seguimientos"
$ws.Range("E13").Value = "Carpeta"

# New row 14: "2.1.1." Formulario de asistencia (Forms)
$ws.Range("A14").Value = "2.1.1."
$ws.Range("B14").Value = "2.1."
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = "Formulario de asistencia"
$ws.Range("E14").Value = "Forms"
$ws.Range("F14").Value = "Formulario de regristo de asistencia durante capacitaciones y ejercicios de seguimiento de proyectos"
$ws.Range("G14").Value = "https://forms.office.com/Pages/ResponsePage.aspx?id=ruyy8CShsk-ezUgjoSOcEF5r9KyA86dEsQa9DwPMXbFUMVRZQkNJSkVOQVpQMExCMTU4RUJLWEUxUy4u"
$ws.Range("G14").Hyperlinks.Add($ws.Range("G14"), "https://forms.office.com/Pages/ResponsePage.aspx?id=ruyy8CShsk-ezUgjoSOcEF5r9KyA86dEsQa9DwPMXbFUMVRZQkNJSkVOQVpQMExCMTU4RUJLWEUxUy4u") | Out-Null

# Restore the hyperlink-cell look (the Hyperlinks.Add call above creates its own
# style) to match the existing hyperlink style used elsewhere in the column.
$ws.Range("G17").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Final selection/view state
$ws.Range("D16").Select()
